$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row 21: "Number of employees / Assets / Turnover" (bold "title" look, like rows 11 & 17) ---
$ws.Cells.Item(21, 2).Value = "Number of employees"
$ws.Cells.Item(21, 3).Value = "Assets (local currency, unless noted otherwise)"
$ws.Cells.Item(21, 4).Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B21:D21").Font.Bold = $true

# --- Row 22: Micro (blank threshold cells, still written so the row exists) ---
$ws.Cells.Item(22, 1).Value = "Micro"
$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = ""

# --- Row 23: Small (blank threshold cells) ---
$ws.Cells.Item(23, 1).Value = "Small"
$ws.Cells.Item(23, 2).Value = ""
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 4).Value = ""

# --- Row 24: Medium + thresholds (overwrite old INSABA caption that used to live here) ---
$ws.Cells.Item(24, 1).Value = "Medium"
$ws.Cells.Item(24, 2).Value = "<10 Manufacturing<br/><5 Other businesses"
$ws.Cells.Item(24, 3).Value = "< N$ 500,000 Manufacturing<br/>< N$ 100,000 Other businesses"
$ws.Cells.Item(24, 4).Value = "< N$ 1,000,000 Manufacturing<br/>< N$ 250,000 Other businesses"
$ws.Cells.Item(24, 1).Font.Bold = $false
$ws.Cells.Item(24, 1).Font.Italic = $false

# --- Row 25: Large + thresholds (overwrite old source caption that used to live here) ---
$ws.Cells.Item(25, 1).Value = "Large"
$ws.Cells.Item(25, 2).Value = "> 10 Manufacturing<br/>> 5 Other businesses"
$ws.Cells.Item(25, 3).Value = "> N$ 500,000 Manufacturing<br/>> N$ 100,000 Other businesses"
$ws.Cells.Item(25, 4).Value = "> N$ 1,000,000 Manufacturing<br/>> N$ 250,000 Other businesses"
$ws.Cells.Item(25, 1).Font.Bold = $false
$ws.Cells.Item(25, 1).Font.Italic = $false

# --- Rows 26-27 used to hold the old "INSABA" source citation block; clear them completely ---
# (content is recreated below at rows 30-33)
$ws.Range("A26:D27").Clear()

# --- Rows 30-33: INSABA source citation block (moved down from the old 24-27 location) ---
$ws.Cells.Item(30, 1).Value = "INSABA"
$ws.Cells.Item(30, 1).Font.Bold = $true

$ws.Cells.Item(31, 1).Value = "Integrated Southern Africa Business Advisory (INSABA), ""Small and Medium Enterprises in Namibia - A brief situational analysis"", 2006, p. 7. Available at http://www.technosol.de/INSABA/Docs/SME%20in%20Namibia-A%20Situational%20Analysis.pdf"
$ws.Cells.Item(31, 1).Font.Italic = $true

$ws.Cells.Item(32, 1).Value = "INSABA"
$ws.Cells.Item(32, 1).Font.Bold = $true

$ws.Cells.Item(33, 1).Value = "Integrated Southern Africa Business Advisory (INSABA), ""Small and Medium Enterprises in Namibia - A brief situational analysis"", 2006, p. 6, 10 and 11. Available at http://www.technosol.de/INSABA/Docs/SME%20in%20Namibia-A%20Situational%20Analysis.pdf"
$ws.Cells.Item(33, 1).Font.Italic = $true
